$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the entire row 699 ("信者" entry); Excel shifts subsequent rows up automatically.
$ws.Rows(699).Delete()
